# Regenerate save_data column G ("K") with corrected strikeout (K) values
# instead of the old "Strike#" values, per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G (header "K")
$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 3
    27 = 3
    28 = 0
    29 = 1
    30 = 1
    31 = 3
    32 = 1
    33 = 1
    34 = 3
    35 = 2
    36 = 1
    37 = 2
    38 = 0
    39 = 3
    40 = 2
    41 = 2
    42 = 3
    43 = 2
    44 = 2
    45 = 3
    46 = 2
    47 = 1
    48 = 2
    49 = 2
    50 = 1
    51 = 2
    52 = 0
    53 = 3
    54 = 0
    55 = 1
    56 = 0
    57 = 3
    58 = 2
    59 = 0
    60 = 2
    61 = 2
    63 = 2
    64 = 1
    65 = 1
    66 = 1
    67 = 1
    68 = 1
    69 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
